# Natmi following Dr Hou advice
# Update the LR-pairs (Il18 -> Il18r1) results table: recomputed statistics for
# existing sending/target cluster pairs and add the missing pairs so every
# sending cluster (FAPs, M1, M2, Neutro, sCs) is reported against every target
# cluster (FAPs, M1, M2, Neutro).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric metric columns (E..T).
$rows = @(
    ,@("FAPs", "Il18", "Il18r1", "FAPs", 3, 1, 7.053021666666666, 21.159065, 0.08011112358180576, 0.08011112358180575, 3, 1, 0.5515783333333334, 1.654735, 0.06037095152961258, 0.06037095152961258, 3.890293935863889, 35.012645422775, 0.004836384758739999, 0.004836384758739999)
    ,@("FAPs", "Il18", "Il18r1", "M1", 3, 1, 7.053021666666666, 21.159065, 0.08011112358180576, 0.08011112358180575, 2, 0.6666666666666666, 0.1679506666666667, 0.503852, 0.01838241450751834, 0.01838241450751834, 1.184559690931111, 10.66103721838, 0.001472635880343781, 0.00147263588034378)
    ,@("FAPs", "Il18", "Il18r1", "M2", 3, 1, 7.053021666666666, 21.159065, 0.08011112358180576, 0.08011112358180575, 3, 1, 0.334542, 1.003626, 0.03661604824933233, 0.03661604824933234, 2.359531974409999, 21.23578776969, 0.002933352766379625, 0.002933352766379625)
    ,@("FAPs", "Il18", "Il18r1", "Neutro", 3, 1, 7.053021666666666, 21.159065, 0.08011112358180576, 0.08011112358180575, 3, 1, 8.082414666666667, 24.247244, 0.8846305857135367, 0.8846305857135368, 57.00544576298444, 513.0490118668599, 0.07086875017634235, 0.07086875017634235)
    ,@("M1", "Il18", "Il18r1", "FAPs", 3, 1, 35.63076933333333, 106.892308, 0.4047089460773642, 0.4047089460773642, 3, 1, 0.5515783333333334, 1.654735, 0.06037095152961258, 0.06037095152961258, 19.65316036426445, 176.87844327838, 0.02443266416723715, 0.02443266416723714)
    ,@("M1", "Il18", "Il18r1", "M1", 3, 1, 35.63076933333333, 106.892308, 0.4047089460773642, 0.4047089460773642, 2, 0.6666666666666666, 0.1679506666666667, 0.503852, 0.01838241450751834, 0.01838241450751834, 5.984211463379555, 53.85790317041599, 0.007439527601694996, 0.007439527601694996)
    ,@("M1", "Il18", "Il18r1", "M2", 3, 1, 35.63076933333333, 106.892308, 0.4047089460773642, 0.4047089460773642, 3, 1, 0.334542, 1.003626, 0.03661604824933233, 0.03661604824933234, 11.919988834312, 107.279899508808, 0.01481884229650521, 0.01481884229650521)
    ,@("M1", "Il18", "Il18r1", "Neutro", 3, 1, 35.63076933333333, 106.892308, 0.4047089460773642, 0.4047089460773642, 3, 1, 8.082414666666667, 24.247244, 0.8846305857135367, 0.8846305857135368, 287.9826526443503, 2591.843873799152, 0.3580179120119268, 0.3580179120119268)
    ,@("M2", "Il18", "Il18r1", "FAPs", 3, 1, 33.64714066666667, 100.941422, 0.3821780750881576, 0.3821780750881575, 3, 1, 0.5515783333333334, 1.654735, 0.06037095152961258, 0.06037095152961258, 18.55903377035222, 167.03130393317, 0.0230724540468278, 0.02307245404682779)
    ,@("M2", "Il18", "Il18r1", "M1", 3, 1, 33.64714066666667, 100.941422, 0.3821780750881576, 0.3821780750881575, 2, 0.6666666666666666, 0.1679506666666667, 0.503852, 0.01838241450751834, 0.01838241450751834, 5.651059706393777, 50.85953735754399, 0.00702535579195598, 0.007025355791955978)
    ,@("M2", "Il18", "Il18r1", "M2", 3, 1, 33.64714066666667, 100.941422, 0.3821780750881576, 0.3821780750881575, 3, 1, 0.334542, 1.003626, 0.03661604824933233, 0.03661604824933234, 11.256381732908, 101.307435596172, 0.01399385083726493, 0.01399385083726493)
    ,@("M2", "Il18", "Il18r1", "Neutro", 3, 1, 33.64714066666667, 100.941422, 0.3821780750881576, 0.3821780750881575, 3, 1, 8.082414666666667, 24.247244, 0.8846305857135367, 0.8846305857135368, 271.9501432156631, 2447.551288940967, 0.3380864144121088, 0.3380864144121088)
    ,@("Neutro", "Il18", "Il18r1", "FAPs", 3, 1, 8.048222, 24.144666, 0.09141501865831142, 0.09141501865831139, 3, 1, 0.5515783333333334, 1.654735, 0.06037095152961258, 0.06037095152961258, 4.439224877056668, 39.95302389351, 0.005518811660499548, 0.005518811660499547)
    ,@("Neutro", "Il18", "Il18r1", "M1", 3, 1, 8.048222, 24.144666, 0.09141501865831142, 0.09141501865831139, 2, 0.6666666666666666, 0.1679506666666667, 0.503852, 0.01838241450751834, 0.01838241450751834, 1.351704250381333, 12.165338253432, 0.001680428765189603, 0.001680428765189603)
    ,@("Neutro", "Il18", "Il18r1", "M2", 3, 1, 8.048222, 24.144666, 0.09141501865831142, 0.09141501865831139, 3, 1, 0.334542, 1.003626, 0.03661604824933233, 0.03661604824933234, 2.692468284324, 24.232214558916, 0.003347256733906346, 0.003347256733906346)
    ,@("Neutro", "Il18", "Il18r1", "Neutro", 3, 1, 8.048222, 24.144666, 0.09141501865831142, 0.09141501865831139, 3, 1, 8.082414666666667, 24.247244, 0.8846305857135367, 0.8846305857135368, 65.04906753338935, 585.441607800504, 0.08086852149871591, 0.0808685214987159)
    ,@("sCs", "Il18", "Il18r1", "FAPs", 3, 1, 3.661325, 10.983975, 0.04158683659436109, 0.04158683659436108, 3, 1, 0.5515783333333334, 1.654735, 0.06037095152961258, 0.06037095152961258, 2.019507541291667, 18.175567871625, 0.002510636896308092, 0.002510636896308092)
    ,@("sCs", "Il18", "Il18r1", "M1", 3, 1, 3.661325, 10.983975, 0.04158683659436109, 0.04158683659436108, 2, 0.6666666666666666, 0.1679506666666667, 0.503852, 0.01838241450751834, 0.01838241450751834, 0.6149219746333333, 5.534297771699999, 0.0007644664683339777, 0.0007644664683339776)
    ,@("sCs", "Il18", "Il18r1", "M2", 3, 1, 3.661325, 10.983975, 0.04158683659436109, 0.04158683659436108, 3, 1, 0.334542, 1.003626, 0.03661604824933233, 0.03661604824933234, 1.22486698815, 11.02380289335, 0.001522745615276225, 0.001522745615276225)
    ,@("sCs", "Il18", "Il18r1", "Neutro", 3, 1, 3.661325, 10.983975, 0.04158683659436109, 0.04158683659436108, 3, 1, 8.082414666666667, 24.247244, 0.8846305857135367, 0.8846305857135368, 29.59234687943333, 266.3311219148999, 0.03678898761444279, 0.03678898761444279)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
